# raven.docx -> "main branch" version
#
# 1) Append a red "(This is a change - Version for main branch)" note to
#    the very first paragraph (after padding the original sentence with
#    two trailing spaces).
# 2) Delete the trailing "...ank God almighty, we are free at last."
#    paragraph at the end of the document.
# 3) Remove a batch of now-unused custom/heading styles left over from
#    the web import (Heading 2/4 + their linked Char styles, and a few
#    podcast/blog leftover styles) that nothing in the body uses.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. First paragraph: pad with two spaces, then add the red note.
#    Anchor on a plain Range just after the final period (position 34)
#    rather than the paragraph's own Range object, so we never touch the
#    trailing paragraph-mark position math.
# ---------------------------------------------------------------------

$note = $d.Range(34, 34)
$note.InsertAfter("  ")

$redStart = $note.End
$note.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$redEnd = $note.End
$d.Range($redStart, $redEnd).Font.Color = 255

$redStart = $note.End
$note.InsertAfter("rsion for main branch")
$redEnd = $note.End
$d.Range($redStart, $redEnd).Font.Color = 255

$redStart = $note.End
$note.InsertAfter(")")
$redEnd = $note.End
$d.Range($redStart, $redEnd).Font.Color = 255

# ---------------------------------------------------------------------
# 2. Drop the last paragraph ("ank God almighty, we are free at last.")
# ---------------------------------------------------------------------

$lastParaIndex = $d.Paragraphs.Count
$d.Paragraphs.Item($lastParaIndex).Range.Delete()

# ---------------------------------------------------------------------
# 3. Prune unused styles. Deleting shifts every later style down by one,
#    so walk from the highest index to the lowest in a single pass and
#    never re-read a style's index/name after the document has been
#    mutated elsewhere in the same pass.
# ---------------------------------------------------------------------

$obsoleteStyleNames = @(
    "heading 2",
    "heading 4",
    "Heading 2 Char",
    "Heading 4 Char",
    "apple-converted-space",
    "Hyperlink",
    "audio-tool",
    "subscribe",
    "subscribe-more-info",
    "generic-title",
    "podcast-tools__subscribe-links"
)

for ($i = $d.Styles.Count; $i -ge 1; $i--) {
    $styleName = $d.Styles.Item($i).NameLocal
    if ($obsoleteStyleNames -contains $styleName) {
        $d.Styles.Item($i).Delete()
    }
}
